$wb = $excel.ActiveWorkbook

# --- Sheet: All Orders ---
$orders = $wb.Worksheets.Item("All Orders")

# Row 8: order status moved from NEW to CANCELLED, and a cancel reason was added.
$orders.Range("H8").Value = "CANCELLED"
$orders.Range("M8").Value = "test order"

# --- Sheet: Daily Summary ---
$summary = $wb.Worksheets.Item("Daily Summary")

# Row 4 corresponds to 2026-01-13: one order was cancelled, reducing revenue/pending by 50.
$summary.Range("D4").Value = 1
$summary.Range("E4").Value = 290
$summary.Range("G4").Value = 290
